# Generate Report for Handoff
#
# Refreshing the localization-status report: the four "Ready for handoff"
# files (04246858, 0d2b9fff, 96e4d007, ac479718) got a new handoff pass,
# so their timestamps move forward and their priority flips from "low" to
# "ht".
#
# - Overview!G4:G7            Latest HO Xliff Generate Date -> 2016-08-29 06:31:52
# - zh-cn!E4:E7 (Priority)    low -> ht
# - zh-cn!H4:H7 (Handoff dt)  2016-08-29 06:31:31 -> 2016-08-29 06:31:47
# - de-de!E4:E7 (Priority)    low -> ht
# - de-de!H4:H7 (Handoff dt)  2016-08-29 06:31:36 -> 2016-08-29 06:31:52

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

foreach ($r in 4..7) {
    $overview.Range("G$r").Value = "2016-08-29 06:31:52"

    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-29 06:31:47"

    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-29 06:31:52"
}
